$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Site resave updated the "hours/mins/secs" sample row with new numbers.
# The source cells are stored as text (numberStoredAsText), so enter the
# values with a leading apostrophe to keep them as text instead of
# letting Excel auto-convert them to numbers, then reset the style back
# to Normal so no stray "quote prefix" formatting is left behind.

$ws.Range("A2").Value = "'127"
$ws.Range("A2").Style = "Normal"

$ws.Range("B2").Value = "'31"
$ws.Range("B2").Style = "Normal"

$ws.Range("C2").Value = "'21"
$ws.Range("C2").Style = "Normal"
